$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 41 (old rows 41-51 shift down to 43-53).
$ws.Rows("41:42").Insert()

# --- New row 41 (Chirimoya, Macroferia Regional de Talca, Maule) ---
$ws.Range("A41").Value = 5
$ws.Range("B41").Value = "Macroferia Regional de Talca"
$ws.Range("C41").Value = "Maule"
$ws.Range("D41").Value = 44508
$ws.Range("E41").Value = 7
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100107
$ws.Range("H41").Value = "Otros"
$ws.Range("I41").Value = 100107002
$ws.Range("J41").Value = "Chirimoya"
$ws.Range("K41").Value = "Cultivar IV Región"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 50
$ws.Range("N41").Value = 25000
$ws.Range("O41").Value = 25000
$ws.Range("P41").Value = 25000
$ws.Range("Q41").Value = '$/bandeja 10 kilos'
$ws.Range("R41").Value = "Provincia de Limarí"
$ws.Range("S41").Value = 2500
$ws.Range("T41").Value = 10

# --- New row 42 (Chirimoya, Macroferia Regional de Talca, Maule) ---
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 44508
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100107
$ws.Range("H42").Value = "Otros"
$ws.Range("I42").Value = 100107002
$ws.Range("J42").Value = "Chirimoya"
$ws.Range("K42").Value = "Cultivar IV Región"
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 60
$ws.Range("N42").Value = 23000
$ws.Range("O42").Value = 23000
$ws.Range("P42").Value = 23000
$ws.Range("Q42").Value = '$/bandeja 10 kilos'
$ws.Range("R42").Value = "Provincia de Limarí"
$ws.Range("S42").Value = 2300
$ws.Range("T42").Value = 10
